# Update numbers on "Comparativa_Final" and "Detalle_Calculos" sheets to the
# latest preelectoral-poll recalculation (adds regional-party handling,
# values provided by the new extraction pipeline).

$wb = $excel.ActiveWorkbook

# --- Sheet 1: Comparativa_Final (A=Partido, B=CIS Oficial, C=Aldabon-Gemini, D=Diferencia)
$ws1 = $wb.Worksheets.Item("Comparativa_Final")

# PP
$ws1.Cells.Item(2, 2).Value = 35.3
$ws1.Cells.Item(2, 3).Value = 39.1
$ws1.Cells.Item(2, 4).Value = 3.8

# PSOE
$ws1.Cells.Item(3, 2).Value = 26.7
$ws1.Cells.Item(3, 4).Value = 0.8

# VOX
$ws1.Cells.Item(4, 2).Value = 15.1
$ws1.Cells.Item(4, 3).Value = 12.9
$ws1.Cells.Item(4, 4).Value = -2.2

# SUMAR
$ws1.Cells.Item(5, 2).Value = 5
$ws1.Cells.Item(5, 3).Value = 2.3
$ws1.Cells.Item(5, 4).Value = -2.7

# --- Sheet 2: Detalle_Calculos
# (A=Partido, B=Real_2023, C=Recuerdo_CIS, D=K_Ponderacion, E=Voto_Simpatia_CIS,
#  F=Ajuste_Fidelidad, G=Final_%)
$ws2 = $wb.Worksheets.Item("Detalle_Calculos")

# PP (row 2)
$ws2.Cells.Item(2, 3).Value = 845.200004372708
$ws2.Cells.Item(2, 4).Value = 0.987
$ws2.Cells.Item(2, 5).Value = 27.8
$ws2.Cells.Item(2, 7).Value = 39.1

# PSOE (row 3)
$ws2.Cells.Item(3, 3).Value = 908.1569883998282
$ws2.Cells.Item(3, 4).Value = 0.88
$ws2.Cells.Item(3, 5).Value = 21.5

# VOX (row 4)
$ws2.Cells.Item(4, 3).Value = 366.9427356363201
$ws2.Cells.Item(4, 4).Value = 0.852
$ws2.Cells.Item(4, 5).Value = 12
$ws2.Cells.Item(4, 7).Value = 12.9

# SUMAR (row 5)
$ws2.Cells.Item(5, 3).Value = 397.34065786201
$ws2.Cells.Item(5, 4).Value = 0.78
$ws2.Cells.Item(5, 5).Value = 2.3
$ws2.Cells.Item(5, 7).Value = 2.3

# ERC (row 6)
$ws2.Cells.Item(6, 3).Value = 0
$ws2.Cells.Item(6, 4).Value = 1
$ws2.Cells.Item(6, 5).Value = 0
$ws2.Cells.Item(6, 7).Value = 0

# JUNTS (row 7)
$ws2.Cells.Item(7, 3).Value = 0
$ws2.Cells.Item(7, 4).Value = 1
$ws2.Cells.Item(7, 5).Value = 0
$ws2.Cells.Item(7, 7).Value = 0

# EH BILDU (row 8)
$ws2.Cells.Item(8, 3).Value = 0
$ws2.Cells.Item(8, 4).Value = 1
$ws2.Cells.Item(8, 5).Value = 0
$ws2.Cells.Item(8, 7).Value = 0

# EAJ-PNV (row 9)
$ws2.Cells.Item(9, 3).Value = 0
$ws2.Cells.Item(9, 4).Value = 1
$ws2.Cells.Item(9, 5).Value = 0
$ws2.Cells.Item(9, 7).Value = 0

# BNG (row 10)
$ws2.Cells.Item(10, 3).Value = 0
$ws2.Cells.Item(10, 4).Value = 1
$ws2.Cells.Item(10, 5).Value = 0
$ws2.Cells.Item(10, 7).Value = 0

# CC (row 11)
$ws2.Cells.Item(11, 3).Value = 0
$ws2.Cells.Item(11, 4).Value = 1
$ws2.Cells.Item(11, 5).Value = 0
$ws2.Cells.Item(11, 7).Value = 0

# UPN (row 12)
$ws2.Cells.Item(12, 3).Value = 0
$ws2.Cells.Item(12, 4).Value = 1
$ws2.Cells.Item(12, 5).Value = 0
$ws2.Cells.Item(12, 7).Value = 0

# PACMA (row 13)
$ws2.Cells.Item(13, 3).Value = 2.73624001264
$ws2.Cells.Item(13, 4).Value = 6.448
